$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1)
$ws.Range("J1").Value = "Wrong number format with dot should be String"
$ws.Range("K1").Value = "Wrong number format with comma should be String"
$ws.Range("L1").Value = "Column contains number should be String"

# Row 2 new data cells
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 123
$ws.Range("L2").Value = 123412341231
$ws.Range("L2").NumberFormat = "General"

# Row 3 new data cells (text values that look like numbers)
$ws.Range("J3").Value = "15024.00.00"
$ws.Range("K3").Value = "15024,00,00"
$ws.Range("L3").Value = "Lorem Ipsum"

# Column widths for the new columns
$ws.Columns.Item(9).ColumnWidth = 73.88095238095238
$ws.Columns.Item(10).ColumnWidth = 41.58333333333333
$ws.Columns.Item(11).ColumnWidth = 47.5
$ws.Columns.Item(12).ColumnWidth = 48.666666666666664

# Selection moves to L1
$ws.Range("L1").Select() | Out-Null
